# Added serial number to identify which table needs to be inserted first.
# Each new "serial number" cell gets a green fill (RGB 0,176,80 -> FF00B050)
# so the order in which the tables below should be read/inserted is obvious.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$green = 5287936   # OLE BGR value for RGB(0,176,80) == FF00B050

$serials = @(
    @{addr="B1";  val=2},
    @{addr="J1";  val=4},
    @{addr="B6";  val=5},
    @{addr="J8";  val=3},
    @{addr="B12"; val=7},
    @{addr="J15"; val=1},
    @{addr="J20"; val=6},
    @{addr="B37"; val=8},
    @{addr="B43"; val=9}
)

foreach ($s in $serials) {
    $cell = $ws.Range($s.addr)
    $cell.Value = $s.val
    $cell.Interior.Color = $green
}

# Update selection to reflect where the user ended up after the edit.
[void]$ws.Range("C22").Select()

Write-Host "done"
